$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-13 (replaces old rows 2-17; target-cluster "ECs" rows removed).
$data = @(
    @("ECs", "Bsg", "Slc16a7", "FAPs", 3, 1, 46.62552400000001, 139.876572, 0.1981356926336952, 0.1981356926336952, 3, 1, 0.450595, 1.351785, 0.02359906163447434, 0.02359906163447434, 21.00922798678, 189.08305188102, 0.004675816422451838, 0.004675816422451837),
    @("ECs", "Bsg", "Slc16a7", "MuSCs", 3, 1, 46.62552400000001, 139.876572, 0.1981356926336952, 0.1981356926336952, 3, 1, 0.66954, 2.00862, 0.03506589226854703, 0.03506589226854703, 31.21765333896001, 280.95888005064, 0.006947804852447105, 0.006947804852447103),
    @("ECs", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 46.62552400000001, 139.876572, 0.1981356926336952, 0.1981356926336952, 3, 1, 17.97363266666667, 53.920898, 0.9413350460969786, 0.9413350460969786, 838.0300412668507, 7542.270371401657, 0.1865120713587963, 0.1865120713587962),
    @("FAPs", "Bsg", "Slc16a7", "FAPs", 3, 1, 103.9192913333333, 311.757874, 0.441606206212991, 0.4416062062129909, 3, 1, 0.450595, 1.351785, 0.02359906163447434, 0.02359906163447434, 46.82551307834333, 421.42961770509, 0.01042149207858676, 0.01042149207858676),
    @("FAPs", "Bsg", "Slc16a7", "MuSCs", 3, 1, 103.9192913333333, 311.757874, 0.441606206212991, 0.4416062062129909, 3, 1, 0.66954, 2.00862, 0.03506589226854703, 0.03506589226854703, 69.57812231932, 626.2031008738801, 0.01548531565218651, 0.0154853156521865),
    @("FAPs", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 103.9192913333333, 311.757874, 0.441606206212991, 0.4416062062129909, 3, 1, 17.97363266666667, 53.920898, 0.9413350460969786, 0.9413350460969786, 1867.80716940565, 16810.26452465085, 0.4156993984822177, 0.4156993984822177),
    @("MuSCs", "Bsg", "Slc16a7", "FAPs", 3, 1, 49.27528633333333, 147.825859, 0.2093958876983056, 0.2093958876983056, 3, 1, 0.450595, 1.351785, 0.02359906163447434, 0.02359906163447434, 22.20319764536833, 199.828778808315, 0.004941546459797782, 0.004941546459797781),
    @("MuSCs", "Bsg", "Slc16a7", "MuSCs", 3, 1, 49.27528633333333, 147.825859, 0.2093958876983056, 0.2093958876983056, 3, 1, 0.66954, 2.00862, 0.03506589226854703, 0.03506589226854703, 32.99177521161999, 296.92597690458, 0.007342653639505558, 0.007342653639505556),
    @("MuSCs", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 49.27528633333333, 147.825859, 0.2093958876983056, 0.2093958876983056, 3, 1, 17.97363266666667, 53.920898, 0.9413350460969786, 0.9413350460969786, 885.6558961001534, 7970.903064901381, 0.1971116875990023, 0.1971116875990022),
    @("Resolving-Mac", "Bsg", "Slc16a7", "FAPs", 3, 1, 35.50107333333333, 106.50322, 0.1508622134550082, 0.1508622134550081, 3, 1, 0.450595, 1.351785, 0.02359906163447434, 0.02359906163447434, 15.99660613863333, 143.9694552477, 0.003560206673637962, 0.003560206673637962),
    @("Resolving-Mac", "Bsg", "Slc16a7", "MuSCs", 3, 1, 35.50107333333333, 106.50322, 0.1508622134550082, 0.1508622134550081, 3, 1, 0.66954, 2.00862, 0.03506589226854703, 0.03506589226854703, 23.7693886396, 213.9244977564, 0.005290118124407863, 0.005290118124407862),
    @("Resolving-Mac", "Bsg", "Slc16a7", "Resolving-Mac", 3, 1, 35.50107333333333, 106.50322, 0.1508622134550082, 0.1508622134550081, 3, 1, 17.97363266666667, 53.920898, 0.9413350460969786, 0.9413350460969786, 638.0832513657289, 5742.74926229156, 0.1420118886569623, 0.1420118886569623),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

# Remove the now-obsolete trailing rows (old rows 14-17, target cluster "ECs" for sender Resolving-Mac).
$ws.Range("A14:T17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
